# "Add Covered states by sf in sf download list"
#
# Insert a new "Covered_States" column (with its {vendor:covered_state}
# merge-field placeholder) right after the existing "State" column (D),
# shifting every column from the old E onward one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; everything from old E..AC moves to F..AD and
# inherits the style of the cells it displaced (matches the header/value
# row styles s="2"/s="3" used throughout row 1 / row 2).
$ws.Columns("E:E").Insert()

# New header (row 1) + merge-field placeholder (row 2) for the inserted column.
$ws.Range("E1").Value = "Covered_States"
$ws.Range("E2").Value = "{vendor:covered_state}"

# Match the new column's width to its left neighbor (State), same as Excel
# does when a column is inserted next to existing, already-sized columns.
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Selection moved off the far-right block (no more topLeftCell scroll, and
# the active cell/selection is now F12 instead of AA13).
$null = $ws.Range("F12").Select()

# Page orientation was switched to Portrait via Page Setup.
$ws.PageSetup.Orientation = 1
